# Auto-generated edit script: update cryptos list values per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.888.25"
$ws.Range("E2").Value = "  +3.06%  "
$ws.Range("D3").Value = "2.469.93"
$ws.Range("E3").Value = "  +1.70%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "489.56"
$ws.Range("E5").Value = "  +2.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.83"
$ws.Range("E6").Value = "  +9.42%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  +2.64%  "
$ws.Range("D9").Value = "2.480.95"
$ws.Range("E9").Value = "  +1.16%  "
$ws.Range("E10").Value = "  +4.16%  "
$ws.Range("E11").Value = "  +4.47%  "
$ws.Range("E12").Value = "  +3.94%  "
$ws.Range("E13").Value = "  +1.60%  "
$ws.Range("D14").Value = "2.907.42"
$ws.Range("E14").Value = "  +2.28%  "
$ws.Range("D15").Value = "57.156.21"
$ws.Range("E15").Value = "  +3.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.01"
$ws.Range("E16").Value = "  +2.69%  "
$ws.Range("E17").Value = "  +2.37%  "
$ws.Range("D18").Value = "2.482.45"
$ws.Range("E18").Value = "  +1.51%  "
$ws.Range("E19").Value = "  +4.78%  "
$ws.Range("E20").Value = "  +3.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "320.64"
$ws.Range("E21").Value = "  +2.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  +0.41%  "
$ws.Range("E23").Value = "  +4.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "58.09"
$ws.Range("E24").Value = "  +1.80%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.163"
$ws.Range("E27").Value = "  +2.09%  "
$ws.Range("D28").Value = "2.596.47"
$ws.Range("E28").Value = "  +2.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.56"
$ws.Range("E29").Value = "  +3.14%  "
$ws.Range("D30").Value = "0.0₃0806"
$ws.Range("E30").Value = "  +5.25%  "
$ws.Range("E31").Value = "  +0.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "151.01"
$ws.Range("E32").Value = "  +1.70%  "
$ws.Range("E33").Value = "  +2.29%  "
$ws.Range("E34").Value = "  +3.43%  "
$ws.Range("E35").Value = "  +1.45%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.15"
$ws.Range("E36").Value = "  +3.17%  "
$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.891"
$ws.Range("E37").Value = "  +6.16%  "
$ws.Range("E38").Value = "  +5.50%  "
$ws.Range("E39").Value = "  +8.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "34.09"
$ws.Range("E40").Value = "  +1.99%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.50"
$ws.Range("E41").Value = "  +3.23%  "
$ws.Range("E42").Value = "  +3.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.995"
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.607"
$ws.Range("E44").Value = "  +1.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0944"
$ws.Range("E45").Value = "  +6.11%  "
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "263.55"
$ws.Range("E46").Value = "  +3.63%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.79"
$ws.Range("E47").Value = "  +3.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.23"
$ws.Range("E48").Value = "  +0.89%  "
$ws.Range("E49").Value = "  +2.88%  "
$ws.Range("E50").Value = "  +4.28%  "
$ws.Range("D51").Value = "1.858.34"
$ws.Range("E51").Value = "  -2.94%  "
